# Update the AssetList worksheet with a fresh set of asset entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the three data rows with the new Asset ID / RFID Tag values.
$ws.Range("A2").Value = "e1264342"
$ws.Range("B2").Value = 45345

$ws.Range("A3").Value = "e1264343"
$ws.Range("B3").Value = 45245345

$ws.Range("A4").Value = "e1264344"
$ws.Range("B4").Value = "waweaser"

# The table now only has 4 rows (was 5) - remove the old trailing row.
$ws.Range("A5:B5").Delete()

# Selection moved to D8 in the saved file.
$ws.Range("D8").Select()
